# Refresh the "cryptos" price/volume snapshot (Price in column D,
# Volume(1h) in column E) for rows 2-51, as produced by the
# GitHub Actions crypto-list updater.
#
# Numeric-looking price strings (e.g. "1.003", "14.00") are written with a
# leading literal apostrophe so Excel keeps them as text (preserving exact
# formatting / trailing zeros) instead of silently recasting them as
# numbers - exactly what typing '1.003 into a cell does natively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.562.33'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '1.802.93'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''328.81'
$ws.Range("E5").Value = '  -2.33%  '
$ws.Range("D6").Value = '''0.9972'
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("D7").Value = '''0.4467'
$ws.Range("E7").Value = '  +4.81%  '
$ws.Range("D8").Value = '''0.3780'
$ws.Range("E8").Value = '  +7.90%  '
$ws.Range("D9").Value = '''44.75'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("D10").Value = '''1.150'
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = '''0.07540'
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").Value = '''22.73'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '''0.9961'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = '''6.322'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").Value = '''7.590'
$ws.Range("E15").Value = '  +3.96%  '
$ws.Range("D16").Value = '1.803.89'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '''0.00001090'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = '''0.06758'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").Value = '''80.94'
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("D20").Value = '''0.9967'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '''17.69'
$ws.Range("E21").Value = '  +2.81%  '
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").Value = '28.624.11'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").Value = '''11.84'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '''2.395'
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").Value = '''20.68'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  -3.63%  '
$ws.Range("D28").Value = '''152.32'
$ws.Range("E28").Value = '  -1.87%  '
$ws.Range("D29").Value = '2.008.20'
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").Value = '''133.32'
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").Value = '''1.264'
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").Value = '''3.935'
$ws.Range("E32").Value = '  -3.38%  '
$ws.Range("D33").Value = '''5.853'
$ws.Range("E33").Value = '  -1.57%  '
$ws.Range("D34").Value = '''0.09355'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").Value = '''0.2277'
$ws.Range("E35").Value = '  +5.09%  '
$ws.Range("D36").Value = '''12.31'
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").Value = '''0.06356'
$ws.Range("E37").Value = '  +1.43%  '
$ws.Range("D38").Value = '''0.02342'
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '''5.225'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.6630'
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = '''1.213'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '''8.154'
$ws.Range("E42").Value = '  +0.98%  '
$ws.Range("D43").Value = '''1.443'
$ws.Range("E43").Value = '  -3.61%  '
$ws.Range("D44").Value = '''0.9964'
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '''14.00'
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("D46").Value = '''0.6104'
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '''3.826'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").Value = '''128.99'
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").Value = '''2.043'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").Value = '''0.07085'
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").Value = '''1.158'
$ws.Range("E51").Value = '  -1.67%  '
